# Update sheet/title for the new "through" date (2022-08-25 -> 2022-08-26)
# and add the new day's carjacking counts for August and the yearly Total.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Through 2022-08-25")

# Rename the worksheet tab
$ws.Name = "Through 2022-08-26"

# Update the header label in I1 ("2022 (through 08-25)" -> "2022 (through 08-26)")
$ws.Range("I1").Value = "2022 (through 08-26)"

# Update August 2022 count (row 9 = August) and yearly Total (row 14)
$ws.Range("I9").Value = 139
$ws.Range("I14").Value = 1110
